{"js": "// Replace each two-digit-by-two-digit multiplication expression in the\n// document's table cells with its updated value. The mapping below is an\n// ordered, 1:1 list of (old, new) text pairs matching the source diff.\nconst replacements = [\n  [\"26\u00d737=\", \"17\u00d782=\"],\n  [\"24\u00d732=\", \"18\u00d796=\"],\n  [\"71\u00d796=\", \"43\u00d723=\"],\n  [\"14\u00d735=\", \"36\u00d768=\"],\n  [\"27\u00d799=\", \"17\u00d780=\"],\n  [\"93\u00d746=\", \"53\u00d728=\"],\n  [\"66\u00d796=\", \"20\u00d741=\"],\n  [\"46\u00d743=\", \"76\u00d733=\"],\n  [\"63\u00d752=\", \"61\u00d726=\"],\n  [\"48\u00d773=\", \"58\u00d793=\"],\n  [\"94\u00d767=\", \"85\u00d773=\"],\n  [\"60\u00d744=\", \"72\u00d716=\"],\n  [\"13\u00d787=\", \"20\u00d721=\"],\n  [\"14\u00d736=\", \"42\u00d759=\"],\n  [\"75\u00d722=\", \"50\u00d732=\"],\n  [\"16\u00d741=\", \"86\u00d773=\"],\n  [\"60\u00d726=\", \"47\u00d724=\"],\n  [\"30\u00d793=\", \"96\u00d719=\"],\n  [\"24\u00d789=\", \"68\u00d794=\"],\n  [\"21\u00d750=\", \"60\u00d742=\"],\n  [\"85\u00d715=\", \"52\u00d771=\"],\n  [\"38\u00d738=\", \"66\u00d732=\"],\n  [\"38\u00d779=\", \"83\u00d778=\"],\n  [\"84\u00d716=\", \"51\u00d742=\"],\n  [\"58\u00d790=\", \"61\u00d777=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  // Each old value is unique in the document, so replace every match found\n  // (expected to be exactly one).\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered, 1:1 list of (old, new) text pairs matching the source diff.\n$pairs = @(\n  @(\"26\u00d737=\", \"17\u00d782=\"),\n  @(\"24\u00d732=\", \"18\u00d796=\"),\n  @(\"71\u00d796=\", \"43\u00d723=\"),\n  @(\"14\u00d735=\", \"36\u00d768=\"),\n  @(\"27\u00d799=\", \"17\u00d780=\"),\n  @(\"93\u00d746=\", \"53\u00d728=\"),\n  @(\"66\u00d796=\", \"20\u00d741=\"),\n  @(\"46\u00d743=\", \"76\u00d733=\"),\n  @(\"63\u00d752=\", \"61\u00d726=\"),\n  @(\"48\u00d773=\", \"58\u00d793=\"),\n  @(\"94\u00d767=\", \"85\u00d773=\"),\n  @(\"60\u00d744=\", \"72\u00d716=\"),\n  @(\"13\u00d787=\", \"20\u00d721=\"),\n  @(\"14\u00d736=\", \"42\u00d759=\"),\n  @(\"75\u00d722=\", \"50\u00d732=\"),\n  @(\"16\u00d741=\", \"86\u00d773=\"),\n  @(\"60\u00d726=\", \"47\u00d724=\"),\n  @(\"30\u00d793=\", \"96\u00d719=\"),\n  @(\"24\u00d789=\", \"68\u00d794=\"),\n  @(\"21\u00d750=\", \"60\u00d742=\"),\n  @(\"85\u00d715=\", \"52\u00d771=\"),\n  @(\"38\u00d738=\", \"66\u00d732=\"),\n  @(\"38\u00d779=\", \"83\u00d778=\"),\n  @(\"84\u00d716=\", \"51\u00d742=\"),\n  @(\"58\u00d790=\", \"61\u00d777=\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($pair[0], $true, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n  if (-not $found) {\n    Write-Host \"NOT FOUND: $($pair[0])\"\n  }\n}\n"}
